# "recent posts and updates"
# The submitted Downey et al. PNAS manuscript row is removed from the
# "citations" sheet, and the dependent TRANSPOSE lookup table on the
# "contributions" sheet is trimmed to match (its "2022Downey" column is
# removed too, since it no longer has a corresponding citation). Finally
# the "contributions" tab is left as the active / selected sheet.

$wb = $excel.ActiveWorkbook

$citations = $wb.Worksheets.Item("citations")
$contributions = $wb.Worksheets.Item("contributions")

# Remove the (submitted) Downey et al. PNAS row - everything below shifts up.
$citations.Rows(2).Delete()

# The contributions sheet carried a per-paper "2022Downey" column (B) driven
# by the now-removed citation; drop it so the TRANSPOSE headers + tallies
# realign with the remaining papers.
$contributions.Range("B1:I1").ClearContents()
$contributions.Columns(2).Delete()
$contributions.Range("B1:H1").FormulaArray = "=TRANSPOSE(citations!E2:E8)"

# Leave the selection where a manual "select row, delete" edit would land.
$citations.Range("A2:XFD2").Select()

# The contributions sheet is the one left showing/active in the saved file.
$contributions.Activate()
$contributions.Range("F13").Select()
